$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2549.6592
$ws.Range("I62").Value = 1355.3077
$ws.Range("J62").Value = 3050.516
$ws.Range("K62").Value = 1355.3077
$ws.Range("L62").Value = 3050.516
$ws.Range("M62").Value = -731.3077000000001
$ws.Range("N62").Value = -4298.516
$ws.Range("H65").Value = 2549.6592
$ws.Range("I65").Value = 1355.3077
$ws.Range("J65").Value = 3050.516
$ws.Range("K65").Value = 6776.538500000001
$ws.Range("L65").Value = 15252.58
$ws.Range("M65").Value = -3656.538500000001
$ws.Range("N65").Value = -21492.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3139041
$ws.Range("I132").Value = 4930907.5
$ws.Range("J132").Value = 3274.5
$ws.Range("K132").Value = 14792722.5
$ws.Range("L132").Value = 9823.5
$ws.Range("M132").Value = -14790192.5
$ws.Range("N132").Value = -14883.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1924.25
$ws.Range("I20").Value = 1626
$ws.Range("J20").Value = 2288.7778
$ws.Range("K20").Value = 1626
$ws.Range("L20").Value = 2288.7778
$ws.Range("M20").Value = -1379
$ws.Range("N20").Value = -2782.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 10300.667
$ws.Range("I33").Value = 8743.714
$ws.Range("J33").Value = 15750
$ws.Range("K33").Value = 8743.714
$ws.Range("L33").Value = 15750
$ws.Range("M33").Value = -8364.714
$ws.Range("N33").Value = -16508

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1424512.8
$ws.Range("I2").Value = 6.3333335
$ws.Range("J2").Value = 2136766
$ws.Range("K2").Value = 38.000001
$ws.Range("L2").Value = 12820596
$ws.Range("M2").Value = 74.999999
$ws.Range("N2").Value = -12820822
$ws.Range("H31").Value = 2900
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("H35").Value = 2292.8572
$ws.Range("J35").Value = 2658.3333
$ws.Range("L35").Value = 7974.999899999999
$ws.Range("N35").Value = -8550.999899999999
$ws.Range("H76").Value = 2600
$ws.Range("I76").Value = 1300
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 3900
$ws.Range("L76").Value = 11700
$ws.Range("M76").Value = -3517
$ws.Range("N76").Value = -12466
$ws.Range("H79").Value = 2600
$ws.Range("I79").Value = 1300
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 3900
$ws.Range("L79").Value = 11700
$ws.Range("M79").Value = -2574
$ws.Range("N79").Value = -14352
$ws.Range("H97").Value = 1840
$ws.Range("I97").Value = 350
$ws.Range("J97").Value = 2026.25
$ws.Range("K97").Value = 1050
$ws.Range("L97").Value = 6078.75
$ws.Range("M97").Value = -554
$ws.Range("N97").Value = -7070.75
$ws.Range("H98").Value = 654.9231
$ws.Range("I98").Value = 740.2
$ws.Range("J98").Value = 601.625
$ws.Range("K98").Value = 2220.6
$ws.Range("L98").Value = 1804.875
$ws.Range("M98").Value = -722.6000000000004
$ws.Range("N98").Value = -4800.875
$ws.Range("H107").Value = 74489.516
$ws.Range("J107").Value = 45893.816
$ws.Range("L107").Value = 137681.448
$ws.Range("N107").Value = -141521.448
$ws.Range("H110").Value = 2400
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H111").Value = 644.25
$ws.Range("I111").Value = 644.25
$ws.Range("K111").Value = 1932.75
$ws.Range("M111").Value = 1134.25
$ws.Range("H131").Value = 1963533.1
$ws.Range("J131").Value = 2327614.5
$ws.Range("L131").Value = 6982843.5
$ws.Range("N131").Value = -6992923.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3090.9092
$ws.Range("J12").Value = 3090.9092
$ws.Range("L12").Value = 3090.9092
$ws.Range("N12").Value = -3370.9092
$ws.Range("H21").Value = 464192.16
$ws.Range("J21").Value = 1501999.5
$ws.Range("L21").Value = 1501999.5
$ws.Range("N21").Value = -1502345.5
$ws.Range("H29").Value = 3000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = ""
$ws.Range("H30").Value = 464192.16
$ws.Range("J30").Value = 1501999.5
$ws.Range("L30").Value = 1501999.5
$ws.Range("N30").Value = -1502209.5
$ws.Range("H62").Value = 29800
$ws.Range("J62").Value = 29800
$ws.Range("L62").Value = 29800
$ws.Range("N62").Value = -31172
$ws.Range("H63").Value = 18998
$ws.Range("I63").Value = 18998
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 18998
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -18312
$ws.Range("N63").Value = ""
$ws.Range("H65").Value = 29800
$ws.Range("J65").Value = 29800
$ws.Range("L65").Value = 89400
$ws.Range("N65").Value = -96264
$ws.Range("H66").Value = 18998
$ws.Range("I66").Value = 18998
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 56994
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -53562
$ws.Range("N66").Value = ""
$ws.Range("H69").Value = 14000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 14000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 14000
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -15498
$ws.Range("H72").Value = 14000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 14000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 42000
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -49488
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""
$ws.Range("H88").Value = 98000
$ws.Range("J88").Value = 98000
$ws.Range("L88").Value = 98000
$ws.Range("N88").Value = -98902
$ws.Range("H91").Value = 98000
$ws.Range("J91").Value = 98000
$ws.Range("L91").Value = 98000
$ws.Range("N91").Value = -101120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 12498.5
$ws.Range("I45").Value = 5000
$ws.Range("K45").Value = 5000
$ws.Range("M45").Value = -4593
$ws.Range("H64").Value = 25716.666
$ws.Range("J64").Value = 25716.666
$ws.Range("L64").Value = 25716.666
$ws.Range("N64").Value = -26166.666
$ws.Range("H67").Value = 25716.666
$ws.Range("J67").Value = 25716.666
$ws.Range("L67").Value = 25716.666
$ws.Range("N67").Value = -27276.666
